$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$shape2 = $s.Shapes.Item(2)
$shape2.Name = "PPRect#2"
$shape2.AutoShapeType = 1

$shape3 = $s.Shapes.Item(3)
$shape3.Name = "PPRect#3"
$shape3.AutoShapeType = 1
